$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 5 (the old 'b3tablocker', 'Ahmed', and 'qasim' rows) -
# only the header row (1) and one data row (2) remain.
$ws.Range("A3:M5").EntireRow.Delete()

# Force text formatting on the numeric/date-looking columns so Excel
# doesn't auto-coerce the typed values to real numbers/dates (the source
# data stores them as plain text), then set the values.
$ws.Range("D2:G2").NumberFormat = "@"
$ws.Range("I2").NumberFormat = "@"

# Overwrite row 2 with the new data (previously the 'qasim' row, with
# an updated dob and an added 'members' privilege).
$ws.Range("A2").Value = "qasim"
$ws.Range("B2").Value = "qasim"
$ws.Range("C2").Value = "Qasim"
$ws.Range("D2").Value = "923432928333"
$ws.Range("E2").Value = "karachi"
$ws.Range("F2").Value = "2025-03-07"
$ws.Range("G2").Value = "20"
$ws.Range("H2").Value = "Male"
$ws.Range("I2").Value = "250000"
$ws.Range("J2").Value = "x"
$ws.Range("K2").Value = "x"
$ws.Range("L2").Value = "members,attendance,payments,packages"
$ws.Range("M2").Value = "trainer"

# Restore the plain/default style on those cells so only the values (not
# the formatting) changed - matches the source which has no style index.
$ws.Range("D2:G2").Style = "Normal"
$ws.Range("I2").Style = "Normal"
